# SetPermissions.xlsx: fill in the "Revise" result for the permission test
# cases and record the Expected/Actual result + Pass/Fail outcome.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (หัวหน้าโครงการ) keeps "Y" but row 3 (สมาชิก) is revised from "N" to "Y".
$ws.Range("A3").Value = "Y"

# New "Expected Result" / "Actual Result" text and "Pass" outcome for both
# test-case rows (these cells were previously empty).
$ws.Range("E2").Style = "Normal"
$ws.Range("E2").Value = "บันทึกข้อมูลสำเร็จ"

$ws.Range("F2").Style = "Normal"
$ws.Range("F2").Value = "Pass"

$ws.Range("E3").Style = "Normal"
$ws.Range("E3").Value = "บันทึกข้อมูลสำเร็จ"

$ws.Range("F3").Style = "Normal"
$ws.Range("F3").Value = "Pass"

# Move the selection to A3.
$ws.Range("A3").Select() | Out-Null

# Column A was effectively relying on the sheet's default width; give it an
# explicit custom width equal to that same default (8.796875 characters is
# the closest value the width-in-characters API can reproduce).
$ws.Columns.Item(1).ColumnWidth = 7.916666666666667

# Column E widens to match column D now that it holds the same
# "บันทึกข้อมูลสำเร็จ" text (closest reproducible value to 13.296875).
$ws.Columns.Item(5).ColumnWidth = 12.416666666666666
